$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "26.740.10"
Set-TextValue $ws.Range("E2") "  +1.06%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.646.57"
Set-TextValue $ws.Range("E3") "  +1.30%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.26%  "

# Row 5
Set-TextValue $ws.Range("D5") "215.94"
Set-TextValue $ws.Range("E5") "  +1.46%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.506"
Set-TextValue $ws.Range("E6") "  +1.53%  "

# Row 7
Set-TextValue $ws.Range("E7") "  +0.15%  "

# Row 8
Set-TextValue $ws.Range("E8") "  +1.84%  "

# Row 9
Set-TextValue $ws.Range("E9") "  +0.66%  "

# Row 10
Set-TextValue $ws.Range("D10") "19.20"
Set-TextValue $ws.Range("E10") "  +2.15%  "

# Row 11
Set-TextValue $ws.Range("E11") "  +0.21%  "

# Row 12
Set-TextValue $ws.Range("D12") "1.876.47"
Set-TextValue $ws.Range("E12") "  +1.27%  "

# Row 13
Set-TextValue $ws.Range("B13") "Polkadot"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "4.18"
Set-TextValue $ws.Range("E13") "  +1.42%  "

# Row 14
Set-TextValue $ws.Range("B14") "WrappedEther"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D14") "1.605.61"
Set-TextValue $ws.Range("E14") "  -0.97%  "

# Row 15
Set-TextValue $ws.Range("E15") "  +1.93%  "

# Row 16
Set-TextValue $ws.Range("D16") "65.47"
Set-TextValue $ws.Range("E16") "  +0.92%  "

# Row 17
Set-TextValue $ws.Range("D17") "26.748.49"
Set-TextValue $ws.Range("E17") "  +0.84%  "

# Row 18
Set-TextValue $ws.Range("E18") "  +0.77%  "

# Row 19
Set-TextValue $ws.Range("D19") "219.22"
Set-TextValue $ws.Range("E19") "  +2.60%  "

# Row 20
Set-TextValue $ws.Range("E20") "  +0.20%  "

# Row 21
Set-TextValue $ws.Range("E21") "  +1.69%  "

# Row 22
Set-TextValue $ws.Range("E22") "  +0.54%  "

# Row 23
Set-TextValue $ws.Range("D23") "2.35"
Set-TextValue $ws.Range("E23") "  +16.97%  "

# Row 24
Set-TextValue $ws.Range("D24") "9.52"
Set-TextValue $ws.Range("E24") "  +2.59%  "

# Row 25
Set-TextValue $ws.Range("D25") "146.10"
Set-TextValue $ws.Range("E25") "  -1.87%  "

# Row 26
Set-TextValue $ws.Range("E26") "  +0.26%  "

# Row 27
Set-TextValue $ws.Range("E27") "  +0.75%  "

# Row 28
Set-TextValue $ws.Range("D28") "7.11"
Set-TextValue $ws.Range("E28") "  +4.21%  "

# Row 29
Set-TextValue $ws.Range("D29") "15.75"
Set-TextValue $ws.Range("E29") "  +1.54%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.0516"
Set-TextValue $ws.Range("E30") "  +1.76%  "

# Row 31
Set-TextValue $ws.Range("E31") "  +1.64%  "

# Row 32
Set-TextValue $ws.Range("E32") "  +0.98%  "

# Row 33
Set-TextValue $ws.Range("E33") "  +2.68%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.281.82"
Set-TextValue $ws.Range("E34") "  +5.11%  "

# Row 35
Set-TextValue $ws.Range("E35") "  +3.79%  "

# Row 36
Set-TextValue $ws.Range("E36") "  +1.95%  "

# Row 37
Set-TextValue $ws.Range("E37") "  +3.35%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.537"
Set-TextValue $ws.Range("E38") "  +6.30%  "

# Row 39
Set-TextValue $ws.Range("E39") "  +3.87%  "

# Row 40
Set-TextValue $ws.Range("E40") "  +0.18%  "

# Row 42
Set-TextValue $ws.Range("E42") "  -0.55%  "

# Row 43
Set-TextValue $ws.Range("E43") "  +2.10%  "

# Row 44
Set-TextValue $ws.Range("D44") "1.787.14"
Set-TextValue $ws.Range("E44") "  +1.33%  "

# Row 45
Set-TextValue $ws.Range("D45") "91.81"
Set-TextValue $ws.Range("E45") "  -1.18%  "

# Row 46
Set-TextValue $ws.Range("D46") "59.75"
Set-TextValue $ws.Range("E46") "  +9.17%  "

# Row 47
Set-TextValue $ws.Range("E47") "  +1.67%  "

# Row 48
Set-TextValue $ws.Range("D48") "0.0516"
Set-TextValue $ws.Range("E48") "  +1.14%  "

# Row 49
Set-TextValue $ws.Range("E49") "  +3.68%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.0969"
Set-TextValue $ws.Range("E50") "  +2.05%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.407"
Set-TextValue $ws.Range("E51") "  +0.12%  "
